# Split the title run "para Oficina Mecânica" into four runs so the
# subtitle becomes "para Oficina Mecânica – SoftGear", matching the
# run-level formatting (rFonts/color/shadow/textOutline) already used by
# the sibling runs in this title paragraph, and flag "SoftGear" with the
# spell-check proofErr markers.

$d = $word.ActiveDocument

# Locate the exact run text we need to split/extend.
$rng = $d.Content
$found = $rng.Find.Execute("para Oficina Mecânica", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the title run 'para Oficina Mecânica'"
}

# Run properties shared by every run on this title line.
$rPr = '<w:rPr>' +
    '<w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi"/>' +
    '<w:color w:val="000000" w:themeColor="text1"/>' +
    '<w14:shadow w14:blurRad="38100" w14:dist="19050" w14:dir="2700000" w14:sx="100000" w14:sy="100000" w14:kx="0" w14:ky="0" w14:algn="tl">' +
        '<w14:schemeClr w14:val="dk1"><w14:alpha w14:val="60000"/></w14:schemeClr>' +
    '</w14:shadow>' +
    '<w14:textOutline w14:w="0" w14:cap="flat" w14:cmpd="sng" w14:algn="ctr">' +
        '<w14:noFill/><w14:prstDash w14:val="solid"/><w14:round/>' +
    '</w14:textOutline>' +
'</w:rPr>'

# Rebuild the found range as four runs: the trimmed "para " run keeps the
# original run's rsidRPr, the three appended runs ("Oficina Mecânica ",
# "– " and "SoftGear") are brand new runs with matching formatting;
# "SoftGear" is wrapped in spellStart/spellEnd proofErr markers.
$body =
    '<w:r w:rsidRPr="00CB2F5F">' + $rPr + '<w:t xml:space="preserve">para </w:t></w:r>' +
    '<w:r>' + $rPr + '<w:t xml:space="preserve">Oficina Mecânica </w:t></w:r>' +
    '<w:r>' + $rPr + '<w:t xml:space="preserve">– </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r>' + $rPr + '<w:t>SoftGear</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>'

$xml = '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
    '<w:body><w:p>' + $body + '</w:p></w:body></w:document>'

$rng.InsertXML($xml)
